# Integrating excel to dataprovider to parameterized data
#
# Adds a new "Sheet4" worksheet at the end of the workbook containing a
# small parameterized-data table, makes it the active/selected sheet, and
# left-aligns its cell contents.

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet.
$lastIndex = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($lastIndex))
$ws.Name = "Sheet4"

# Column A: Hello1 / Hello2 / Hello3
$ws.Range("A2").Value = "Hello1"
$ws.Range("A3").Value = "Hello2 "
$ws.Range("A4").Value = "Hello3 "

# Column B (rows 4 then 3): Text3 / Text2
$ws.Range("B4").Value = "Text3 "
$ws.Range("B3").Value = "Text2 "

# Header row: Param1 / Param2 / Param3
$ws.Range("A1").Value = "Param1"
$ws.Range("B1").Value = "Param2"
$ws.Range("C1").Value = "Param3"

# Remaining cell
$ws.Range("B2").Value = "Text1"

# Column C numeric values
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 3

# Left-align all populated cells
$ws.Range("A1:C4").HorizontalAlignment = -4131

# Column widths
$ws.Columns.Item(1).ColumnWidth = 7.666666666666666
$ws.Columns.Item(2).ColumnWidth = 7.5
$ws.Columns.Item(3).ColumnWidth = 7.166666666666666

# Selection / active cell on the new sheet
$ws.Range("M8").Select()
